# Update NATMI LR-pair stats (Hras-Agtr1a) with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values per row, keyed by column letter -> value.
# Only columns whose computed statistic actually changed with the new TPM
# input are included; identifiers (A-F) are left untouched.
$updates = @{
    2  = @{
        G = 6.400771666666667;  H = 19.202315;            I = 0.4226371084084476
        J = 0.4226371084084476; K = 3;                     L = 1
        M = 1.137783333333333;  N = 3.41335
        O = 0.02543039699931523; P = 0.02543039699931523
        Q = 7.282691322805556;  R = 65.54422190525001
        S = 0.01074782945346945; T = 0.01074782945346945
    }
    3  = @{
        G = 6.400771666666667;  H = 19.202315;            I = 0.4226371084084476
        J = 0.4226371084084476
        O = 0.3452795715412271; P = 0.345279571541227
        Q = 98.88027071197612;  R = 889.9224364077851
        S = 0.1459279597086919; T = 0.1459279597086919
    }
    4  = @{
        G = 6.400771666666667;  H = 19.202315;            I = 0.4226371084084476
        J = 0.4226371084084476
        M = 28.15511333333333;  N = 84.46534
        O = 0.6292900314594577; P = 0.6292900314594577
        Q = 180.2144516957889;  R = 1621.9300652621
        S = 0.2659613192462862; T = 0.2659613192462862
    }
    5  = @{
        I = 0.3497933355610079; J = 0.3497933355610079
        K = 3;                  L = 1
        M = 1.137783333333333;  N = 3.41335
        O = 0.02543039699931523; P = 0.02543039699931523
        Q = 6.027480405727778;  R = 54.24732365155001
        S = 0.00889538339103112; T = 0.00889538339103112
    }
    6  = @{
        I = 0.3497933355610079; J = 0.3497933355610079
        O = 0.3452795715412271; P = 0.345279571541227
        S = 0.1207764930304815; T = 0.1207764930304815
    }
    7  = @{
        I = 0.3497933355610079; J = 0.3497933355610079
        M = 28.15511333333333;  N = 84.46534
        O = 0.6292900314594577; P = 0.6292900314594577
        Q = 149.1535241956244;  R = 1342.38171776062
        S = 0.2201214591394954; T = 0.2201214591394954
    }
    8  = @{
        G = 3.446504666666667;  H = 10.339514;            I = 0.2275695560305444
        J = 0.2275695560305443; K = 3;                     L = 1
        M = 1.137783333333333;  N = 3.41335
        O = 0.02543039699931523; P = 0.02543039699931523
        Q = 3.92137556798889;   R = 35.29238011190001
        S = 0.005787184154814653; T = 0.005787184154814652
    }
    9  = @{
        G = 3.446504666666667;  H = 10.339514;            I = 0.2275695560305444
        J = 0.2275695560305443
        O = 0.3452795715412271; P = 0.345279571541227
        Q = 53.24222331267178;  R = 479.180009814046
        S = 0.07857511880205362; T = 0.0785751188020536
    }
    10 = @{
        G = 3.446504666666667;  H = 10.339514;            I = 0.2275695560305444
        J = 0.2275695560305443
        M = 28.15511333333333;  N = 84.46534
        O = 0.6292900314594577; P = 0.6292900314594577
        Q = 97.03672949386224;  R = 873.3305654447601
        S = 0.1432072530736761; T = 0.1432072530736761
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
